$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D2").Value = "SKIP"
$ws.Range("D3").Value = "SKIP"
$ws.Range("D4").Value = "SKIP"
